# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 618893
$ws.Range("C4").Value = 5007
$ws.Range("D4").Value = 47069
$ws.Range("E4").Value = 544738
$ws.Range("G4").Value = 1039
$ws.Range("H4").Value = 27086

# --- Row 8: Alemania ---
$ws.Range("B8").Value = 133154
$ws.Range("C8").Value = 944
$ws.Range("E8").Value = 56962

# --- Row 16: Suiza ---
$ws.Range("D16").Value = 15400
$ws.Range("E16").Value = 9710

# --- Row 29: Ecuador ---
$ws.Range("B29").Value = 7858
$ws.Range("C29").Value = 255
$ws.Range("D29").Value = 780
$ws.Range("E29").Value = 6690
$ws.Range("F29").Value = 135
$ws.Range("G29").Value = 19
$ws.Range("H29").Value = 388

# --- Row 50: Luxemburgo ---
$ws.Range("B50").Value = 3373
$ws.Range("C50").Value = 66
$ws.Range("D50").Value = 526
$ws.Range("E50").Value = 2778
$ws.Range("F50").Value = 33
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 69

# --- Rows 54 & 55: swap Argentina and Sudafrica, with updated stats ---
# Row 54 becomes Sudafrica (new, updated numbers)
$ws.Range("A54").Value = "Sudafrica"
$ws.Range("B54").Value = 2506
$ws.Range("C54").Value = 91
$ws.Range("D54").Value = 410
$ws.Range("E54").Value = 2062
$ws.Range("F54").Value = 7
$ws.Range("G54").Value = 7
$ws.Range("H54").Value = 34

# Row 55 becomes Argentina (keeps prior stats)
$ws.Range("A55").Value = "Argentina"
$ws.Range("B55").Value = 2443
$ws.Range("C55").Value = 166
$ws.Range("D55").Value = 596
$ws.Range("E55").Value = 1738
$ws.Range("F55").Value = 117
$ws.Range("G55").Value = 7
$ws.Range("H55").Value = 109

# --- Row 59: Moldavia ---
$ws.Range("E59").Value = 1832
$ws.Range("G59").Value = 5
$ws.Range("H59").Value = 46

# --- Row 69: Kazajistan ---
$ws.Range("B69").Value = 1295
$ws.Range("C69").Value = 63
$ws.Range("E69").Value = 1039

# --- Updated timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 18:52"
